$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A3"   = -22.0597
    "A14"  = -21.7799
    "B15"  = 5.113099999999998
    "A16"  = -22.08190000000002
    "A21"  = -21.8387
    "B21"  = 5.6298
    "B22"  = 10.0944
    "A23"  = -20.05329999999998
    "B24"  = 5.852000000000001
    "A25"  = -21.77439999999999
    "A26"  = -21.06109999999996
    "B27"  = 6.748200000000005
    "B28"  = 5.499700000000002
    "A29"  = -21.41209999999998
    "B36"  = 9.161999999999999
    "B39"  = 9.015600000000001
    "A40"  = -20.113
    "B45"  = 5.332000000000006
    "B48"  = 7.482200000000003
    "B49"  = 5.387099999999996
    "B52"  = 5.446999999999997
    "A53"  = -21.7427
    "B53"  = 5.2485
    "B54"  = 4.797800000000001
    "A57"  = -22.57360000000003
    "B57"  = 4.563299999999996
    "A59"  = -22.52800000000001
    "A65"  = -21.79639999999998
    "A69"  = -21.6449
    "B70"  = 7.965500000000004
    "B71"  = 4.639599999999996
    "A79"  = -20.4023
    "A83"  = -21.67919999999999
    "B86"  = 4.964700000000001
    "B87"  = 5.498999999999994
    "B89"  = 4.815899999999997
    "A91"  = -20.69489999999998
    "A93"  = -21.30060000000002
    "A100" = -22.20130000000001
    "B101" = 4.9384
    "A103" = -21.7274
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
